{"js": "// The edit merges two runs inside one paragraph (\"M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i\n// li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3. \" + \"K\u1ebft lu\u1eadn n\u00e0y ch\u1ec9 ra t\u1ea7m quan tr\u1ecdng ...\")\n// into a single run holding completely new / expanded Vietnamese text\n// about low-carbohydrate diets, while keeping the existing run formatting\n// (w:sz 24).\n\nconst newText =\n  \"Ch\u1ebf \u0111\u1ed9 \u0103n \u00edt carbohydrate ng\u00e0y c\u00e0ng tr\u1edf n\u00ean ph\u1ed5 bi\u1ebfn. Nh\u1eefng ng\u01b0\u1eddi \u1ee7ng \" +\n  \"h\u1ed9 cho r\u1eb1ng ch\u00fang c\u00f3 hi\u1ec7u qu\u1ea3 gi\u1ea3m c\u00e2n \u0111\u00e1ng k\u1ec3 h\u01a1n c\u00e1c ch\u1ebf \u0111\u1ed9 kh\u00e1c v\u00e0 \" +\n  \"mang l\u1ea1i nh\u1eefng l\u1ee3i \u00edch s\u1ee9c kh\u1ecfe kh\u00e1c nh\u01b0 h\u1ea1 huy\u1ebft \u00e1p v\u00e0 c\u1ea3i thi\u1ec7n \" +\n  \"n\u1ed3ng \u0111\u1ed9 cholesterol; tuy nhi\u00ean, m\u1ed9t s\u1ed1 b\u00e1c s\u0129 tin r\u1eb1ng nh\u1eefng b\u1eefa \u0103n \" +\n  \"ki\u00eang n\u00e0y ti\u1ec1m \u1ea9n nh\u1eefng r\u1ee7i ro v\u1ec1 l\u00e2u d\u00e0i. M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i \" +\n  \"li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3 ch\u1ec9 ra r\u1eb1ng th\u1ef1c \u0111\u01a1n low-carbohydrate r\u1ea5t h\u1eefu \" +\n  \"hi\u1ec7u trong \u1ec7c gi\u1ea3m b\u00e9o ng\u1eafn h\u1ea1n nh\u01b0ng t\u00e1c d\u1ee5ng l\u00e2u b\u1ec1n c\u1ee7a ch\u00fang \" +\n  \"kh\u00f4ng l\u1edbn h\u01a1n nhi\u1ec1u so v\u1edbi c\u00e1c k\u1ebf ho\u1ea1ch \u0103n u\u1ed1ng th\u00f4ng th\u01b0\u1eddng kh\u00e1c. \" +\n  \"\u1ea2nh \";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that currently holds the two runs being merged\n// (identified by the start of the first run's text) instead of assuming\n// a fixed index, so the script is resilient to minor structural changes.\nconst marker = \"M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the paragraph to update.\");\n}\n\n// Replacing the whole paragraph range with new text collapses the\n// paragraph's existing runs into one run and keeps that run's\n// formatting (matches the diff: both <w:r> merge into a single\n// <w:r> with the original <w:rPr><w:sz w:val=\"24\"/></w:rPr>).\ntarget.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The edit merges two runs inside one paragraph (\"M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i\n# li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3. \" + \"K\u1ebft lu\u1eadn n\u00e0y ch\u1ec9 ra t\u1ea7m quan tr\u1ecdng ...\")\n# into a single run holding completely new / expanded Vietnamese text\n# about low-carbohydrate diets, while keeping the existing run formatting\n# (w:sz 24).\n\n$d = $word.ActiveDocument\n\n$newText = \"Ch\u1ebf \u0111\u1ed9 \u0103n \u00edt carbohydrate ng\u00e0y c\u00e0ng tr\u1edf n\u00ean ph\u1ed5 bi\u1ebfn. Nh\u1eefng ng\u01b0\u1eddi \u1ee7ng h\u1ed9 cho r\u1eb1ng ch\u00fang c\u00f3 hi\u1ec7u qu\u1ea3 gi\u1ea3m c\u00e2n \u0111\u00e1ng k\u1ec3 h\u01a1n c\u00e1c ch\u1ebf \u0111\u1ed9 kh\u00e1c v\u00e0 mang l\u1ea1i nh\u1eefng l\u1ee3i \u00edch s\u1ee9c kh\u1ecfe kh\u00e1c nh\u01b0 h\u1ea1 huy\u1ebft \u00e1p v\u00e0 c\u1ea3i thi\u1ec7n n\u1ed3ng \u0111\u1ed9 cholesterol; tuy nhi\u00ean, m\u1ed9t s\u1ed1 b\u00e1c s\u0129 tin r\u1eb1ng nh\u1eefng b\u1eefa \u0103n ki\u00eang n\u00e0y ti\u1ec1m \u1ea9n nh\u1eefng r\u1ee7i ro v\u1ec1 l\u00e2u d\u00e0i. M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3 ch\u1ec9 ra r\u1eb1ng th\u1ef1c \u0111\u01a1n low-carbohydrate r\u1ea5t h\u1eefu hi\u1ec7u trong \u1ec7c gi\u1ea3m b\u00e9o ng\u1eafn h\u1ea1n nh\u01b0ng t\u00e1c d\u1ee5ng l\u00e2u b\u1ec1n c\u1ee7a ch\u00fang kh\u00f4ng l\u1edbn h\u01a1n nhi\u1ec1u so v\u1edbi c\u00e1c k\u1ebf ho\u1ea1ch \u0103n u\u1ed1ng th\u00f4ng th\u01b0\u1eddng kh\u00e1c. \u1ea2nh \"\n\n# Locate the paragraph that currently holds the two runs being merged\n# (identified by the start of the first run's text) instead of assuming\n# a fixed index, so the script is resilient to minor structural changes.\n$marker = \"M\u1ed9t \u0111\u00e1nh gi\u00e1 v\u1ec1 c\u00e1c t\u00e0i li\u1ec7u nghi\u00ean c\u1ee9u hi\u1ec7n c\u00f3.\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($marker)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the paragraph to update.\"\n}\n\n# Paragraph.Range includes the trailing paragraph mark, so build a range\n# that stops one character short of it; assigning .Text on that range\n# collapses the paragraph's existing runs into a single run and keeps\n# that run's formatting (matches the diff: both <w:r> merge into one\n# <w:r> with the original <w:rPr><w:sz w:val=\"24\"/></w:rPr>).\n$start = $target.Range.Start\n$end = $target.Range.End\n$r = $d.Range($start, $end - 1)\n$r.Text = $newText\n"}
